# Apply edits described by the diff:
#  - "About" sheet: update source + notes, clear the now-unused rows/hyperlink
#  - "WMITR" sheet: update the marginal income tax rate value

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsWmitr = $wb.Worksheets.Item("WMITR")

# --- About sheet -----------------------------------------------------

# Source changes from CBO to Census/NerdWallet
$wsAbout.Range("B3").Value = "Census and NerdWallet"

# The year (B4) and the old headline/hyperlink rows (B5:B7) are no longer used
$wsAbout.Range("B4").ClearContents()
$wsAbout.Range("B5").ClearContents()

# Remove the old hyperlink on B6 and clear its text
if ($wsAbout.Hyperlinks.Count -gt 0) {
    $wsAbout.Hyperlinks.Delete()
}
$wsAbout.Range("B6").ClearContents()
$wsAbout.Range("B7").ClearContents()

# Notes text replaced with a single consolidated note, old A11 removed
$wsAbout.Range("A10").Value = "California has a median household income of `$78,672, which puts it in the 6% marginal income tax bracket"
$wsAbout.Range("A11").ClearContents()

# --- WMITR sheet -------------------------------------------------------

$wsWmitr.Range("B2").Value = 0.06

# Match the author's final on-screen selection in each sheet
[void]$wsAbout.Range("A11").Select()
[void]$wsWmitr.Range("B3").Select()
[void]$wsAbout.Select()

$wb.Save()
